$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Fix the "Return2Performande" -> "Return2Performance" typo (header S1).
# ---------------------------------------------------------------------------
$ws.Range("S1").Value2 = "Return2Performance"

# ---------------------------------------------------------------------------
# 2) Matrix edit: show individual tiles for "casussen"/"kennis" rows instead
#    of a single counter cell - move the "v" mark for row 4 from the
#    "Schriftelijk product" column (L) to the "Performance" column (O), and
#    add a new "v" mark for row 3 in the "Schriftelijk product" column (L).
# ---------------------------------------------------------------------------
$ws.Range("L4").Clear() | Out-Null
$ws.Range("O4").Value2 = "v"
$ws.Range("L3").Value2 = "v"

# ---------------------------------------------------------------------------
# 3) Per-column best-fit widths (was: a handful of grouped ranges sharing one
#    width each). Re-create the new, per-column layout.
# ---------------------------------------------------------------------------
$widths = @{
     1 = 26.296875
     2 = 26.09765625
     3 = 7.3984375
     4 = 9.796875
     5 = 7.3984375
     6 = 9.69921875
     7 = 9.59765625
     8 = 9.59765625
     9 = 9.5
    10 = 11.19921875
    11 = 6.796875
    12 = 7.69921875
    13 = 7.09765625
    14 = 7.19921875
    15 = 7.59765625
    16 = 7.69921875
    17 = 6.19921875
    18 = 5
    19 = 5.59765625
    20 = 5.3984375
    21 = 5.59765625
    22 = 5.59765625
    23 = 5.59765625
    24 = 5.5
    25 = 5.5
    26 = 5.09765625
    27 = 5.09765625
}
foreach ($col in $widths.Keys) {
    # OOXML <col width> = ColumnWidth (chars) + 5/6, so back it out here.
    $ws.Columns.Item($col).ColumnWidth = $widths[$col] - (5/6)
}

# ---------------------------------------------------------------------------
# 4) View state: freeze panes stay at B2 (1 row / 1 column), but the visible
#    (bottom-right) pane now scrolls to K2, with V3 selected.
# ---------------------------------------------------------------------------
$ws.Range("V3").Select() | Out-Null
